# Updated cryptos list with GitHub Actions
# This script applies the latest price / 1h-volume snapshot to the
# "cryptos" worksheet. Price cells in column D and volume-percentage
# cells in column E are stored as text in the workbook, so any value
# that Excel would otherwise auto-parse as a number is written with a
# leading apostrophe to force it to remain text (matching the original
# t="inlineStr"/shared-string text representation instead of being
# silently converted to a numeric cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.132.47"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").Value = "3.136.16"
$ws.Range("E3").Value = "  +3.30%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'580.75"
$ws.Range("E5").Value = "  +0.46%  "
$ws.Range("D6").Value = "'174.69"
$ws.Range("E6").Value = "  +3.70%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.132.48"
$ws.Range("E8").Value = "  +3.27%  "
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("D10").Value = "'6.50"
$ws.Range("E10").Value = "  -2.98%  "
$ws.Range("E11").Value = "  +1.64%  "
$ws.Range("D12").Value = "'0.485"
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("D13").Value = "'0.0000250"
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("D14").Value = "'37.27"
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("D16").Value = "3.654.94"
$ws.Range("E16").Value = "  +3.33%  "
$ws.Range("D17").Value = "67.147.62"
$ws.Range("E17").Value = "  +1.19%  "
$ws.Range("D18").Value = "'7.19"
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("D19").Value = "3.136.43"
$ws.Range("E19").Value = "  +3.49%  "
$ws.Range("D21").Value = "'486.76"
$ws.Range("E21").Value = "  +3.99%  "
$ws.Range("D22").Value = "'0.718"
$ws.Range("E22").Value = "  +1.38%  "
$ws.Range("D23").Value = "'7.65"
$ws.Range("E23").Value = "  +3.52%  "
$ws.Range("D24").Value = "'84.42"
$ws.Range("E24").Value = "  +1.55%  "
$ws.Range("D25").Value = "'13.25"
$ws.Range("E25").Value = "  +3.41%  "
$ws.Range("E26").Value = "  +2.97%  "
$ws.Range("E27").Value = "  +0.40%  "
$ws.Range("D29").Value = "'7.99"
$ws.Range("E29").Value = "  -2.93%  "
$ws.Range("E30").Value = "  -1.98%  "
$ws.Range("D31").Value = "'2.69"
$ws.Range("E31").Value = "  +1.70%  "
$ws.Range("D32").Value = "'28.88"
$ws.Range("E32").Value = "  +1.88%  "
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("D34").Value = "'0.114"
$ws.Range("E34").Value = "  -3.52%  "
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").Value = "'5.95"
$ws.Range("E36").Value = "  +1.18%  "
$ws.Range("D37").Value = "'0.991"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").Value = "'47.33"
$ws.Range("E38").Value = "  -2.30%  "
$ws.Range("E39").Value = "  +2.31%  "
$ws.Range("D40").Value = "'50.18"
$ws.Range("E41").Value = "  +0.53%  "
$ws.Range("D42").Value = "'0.123"
$ws.Range("E42").Value = "  +1.89%  "
$ws.Range("D43").Value = "'8.68"
$ws.Range("E44").Value = "  -1.28%  "
$ws.Range("D45").Value = "2.859.72"
$ws.Range("E45").Value = "  +5.46%  "
$ws.Range("D46").Value = "'385.82"
$ws.Range("E46").Value = "  +1.59%  "
$ws.Range("D47").Value = "'0.0358"
$ws.Range("E47").Value = "  -0.73%  "
$ws.Range("D48").Value = "'136.32"
$ws.Range("E48").Value = "  +1.03%  "
$ws.Range("D50").Value = "'25.00"
$ws.Range("E50").Value = "  +1.87%  "
$ws.Range("D51").Value = "'2.22"
$ws.Range("E51").Value = "  -0.80%  "
